# Fix formatting issues introduced when scraping floating point numbers:
# 1) A handful of "Razon social" entries used a comma as a name separator,
#    which is ambiguous with the decimal separator used elsewhere in the
#    sheet. Replace those separating commas with periods.
# 2) The "Importe" column was scraped using Argentine number formatting
#    (period as thousands separator, comma as decimal separator) but the
#    values are stored as text. Re-write them using a plain decimal-point
#    format (no thousands separators).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Razon social fixes (comma -> period as separator) ---
$ws.Range("E85").Value = "BOFFELLI. MARIA INES"
$ws.Range("E96").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E140").Value = "GIMENEZ ANIBAL. FALISTOCCO MARISA DANIELA SH"
$ws.Range("E171").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E188").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- 2) Importe column fixes (Argentine "1.234,56" -> "1234.56") ---
$amounts = @(
  "2|52000.00",
  "3|119060.00",
  "4|7200.00",
  "5|64050.00",
  "6|498000.00",
  "7|352800.00",
  "8|498000.00",
  "9|746000.00",
  "10|498000.00",
  "11|735000.00",
  "12|40461.00",
  "13|1450.00",
  "14|889.90",
  "15|105.00",
  "16|10150.00",
  "17|75.00",
  "18|34936.08",
  "19|1204.50",
  "20|17548.75",
  "21|320.00",
  "22|6640.00",
  "23|917178.00",
  "24|221358.64",
  "25|322608.00",
  "26|95259.70",
  "27|81500.00",
  "28|386600.00",
  "29|13050.00",
  "30|18047.60",
  "31|40075.00",
  "32|69932.00",
  "33|16047.60",
  "34|126745.09",
  "35|735.00",
  "36|120800.00",
  "37|39857.99",
  "38|2500.00",
  "39|13000.00",
  "40|7900.00",
  "41|234.00",
  "42|4800.00",
  "43|112000.00",
  "44|600.00",
  "45|533.00",
  "46|6850.00",
  "47|23718.00",
  "48|400.00",
  "49|429439.91",
  "50|38028.39",
  "51|69993.92",
  "52|1322481.60",
  "53|17629.00",
  "54|96206.00",
  "55|3274600.00",
  "56|9963.14",
  "57|550.60",
  "58|1576.61",
  "59|3675.00",
  "60|7710.00",
  "61|35240.86",
  "62|427.80",
  "63|10617.09",
  "64|581405.00",
  "65|4199.28",
  "66|17088.95",
  "67|4525.26",
  "68|13730.00",
  "69|32385.00",
  "70|2771.00",
  "71|2849.85",
  "72|1958.00",
  "73|100.00",
  "74|51456.00",
  "75|8360.00",
  "76|36299.79",
  "77|3300.00",
  "78|27740.00",
  "79|11710.00",
  "80|177.60",
  "81|2944.00",
  "82|11400.00",
  "83|200.00",
  "84|13900.00",
  "85|29400.00",
  "86|1799800.00",
  "87|32000.00",
  "88|87500.00",
  "89|5850.00",
  "90|151500.00",
  "91|24400.00",
  "92|32500.00",
  "93|91255.00",
  "94|660.00",
  "95|1080.00",
  "96|6260.00",
  "97|125830.00",
  "98|3300.00",
  "99|1835.18",
  "100|15619.73",
  "101|300.00",
  "102|5861.00",
  "103|6914.75",
  "104|13900.00",
  "105|730.00",
  "106|844.24",
  "107|400.00",
  "108|6355.05",
  "109|2432.00",
  "110|180.00",
  "111|26680.00",
  "112|38629.50",
  "113|115.00",
  "114|7364.30",
  "115|3950.00",
  "116|102944.00",
  "117|259351.00",
  "118|315.00",
  "119|4297.00",
  "120|28745.00",
  "121|10286.00",
  "122|28557.27",
  "123|2739.00",
  "124|80.00",
  "125|11855.00",
  "126|2484.00",
  "127|750.00",
  "128|11200.00",
  "129|86905.83",
  "130|68789.00",
  "131|20130.00",
  "132|5070.00",
  "133|327046.00",
  "134|40000.00",
  "135|131600.00",
  "136|540000.00",
  "137|14400.00",
  "138|21120.00",
  "139|30647.00",
  "140|187300.00",
  "141|6845.00",
  "142|40132.00",
  "143|65400.00",
  "144|10000.00",
  "145|12000.00",
  "146|10000.00",
  "147|14000.00",
  "148|37000.00",
  "149|12000.00",
  "150|6000.00",
  "151|7000.00",
  "152|10500.00",
  "153|24000.00",
  "154|9500.00",
  "155|33000.00",
  "156|45200.00",
  "157|156000.00",
  "158|53000.00",
  "159|22800.00",
  "160|7610.20",
  "161|19272.00",
  "162|1920.76",
  "163|40390.00",
  "164|951.00",
  "165|3960.00",
  "166|1550.00",
  "167|8000.00",
  "168|9000.00",
  "169|52166.50",
  "170|10000.00",
  "171|8570.00",
  "172|9800.00",
  "173|47035.12",
  "174|9200.00",
  "175|59539.60",
  "176|12700.00",
  "177|213.63",
  "178|21605.00",
  "179|2880.00",
  "180|6975.00",
  "181|1930.00",
  "182|18530.00",
  "183|181550.00",
  "184|28800.00",
  "185|1200.00",
  "186|490.00",
  "187|3604.28",
  "188|1400.00",
  "189|1320.00",
  "190|7500.00",
  "191|7900.00",
  "192|2796.50",
  "193|2000.00",
  "194|287.91",
  "195|1930.00",
  "196|455.70",
  "197|29210.00",
  "198|900.00",
  "199|760.00",
  "200|11599.07",
  "201|9590.88",
  "202|80000.00",
  "203|40000.00",
  "204|40000.00",
  "205|40000.00",
  "206|80000.00",
  "207|40000.00",
  "208|55000.00",
  "209|40000.00",
  "210|40000.00",
  "211|80000.00",
  "212|80000.00",
  "213|15000.00",
  "214|23754.85",
  "215|6840918.99",
  "216|12945.10",
  "217|695000.00",
  "218|5470727.58",
  "219|9000.00",
  "220|250000.00",
  "221|270250.00",
  "222|273700.00",
  "223|253200.00",
  "224|250000.00",
  "225|250000.00",
  "226|480350.00",
  "227|250000.00",
  "228|603450.00",
  "229|557000.00",
  "230|304900.00",
  "231|250000.00",
  "232|250000.00",
  "233|500000.00",
  "234|396300.00",
  "235|494800.00",
  "236|711150.00",
  "237|470000.00",
  "238|725400.00",
  "239|500000.00",
  "240|258800.00",
  "241|13200.00",
  "242|8000.00",
  "243|1432750.25",
  "244|1070000.00",
  "245|309100.00",
  "246|128106.83",
  "247|118000.00",
  "248|10500.00",
  "249|9231.75",
  "250|106080.00",
  "251|9800.00",
  "252|45100.00",
  "253|198800.00",
  "254|3000.00",
  "255|6800.00",
  "256|3000.00",
  "257|4560000.00",
  "258|29635.00",
  "259|4200.00"
)

foreach ($entry in $amounts) {
  $parts = $entry.Split("|")
  $row = $parts[0]
  $newValue = $parts[1]
  $cell = $ws.Cells.Item([int]$row, 8)

  # Preserve the cell's existing style: force a temporary text format so the
  # numeric-looking string isn't auto-converted into a real number (which
  # would drop the trailing zeros / meaningfully change the stored type),
  # then restore the original style once the text value is set.
  $origStyle = $cell.Style
  $cell.NumberFormat = "@"
  $cell.Value = $newValue
  $cell.Style = $origStyle
}
